$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest day's profit figure (run on 2025-09-21) as a new row
# directly below the existing data, following the same layout used for
# every other row: column A holds the date as literal text, column B
# holds the profit as a plain number.
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# Force column A to text formatting before writing so the date string is
# not auto-converted into a date serial number, then restore the default
# (unstyled) look to match the rest of the column.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "09/21/2025"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 15252.04
